$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Load factor" header in D1, matching the existing header style ---
$ws.Range("D1").Value = "Load factor"
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)   # xlPasteFormats

# --- Load factor values for 2005-2019 ---
$loadFactors = @(
    0.776,
    0.792,
    0.799,
    0.795,
    0.804,
    0.821,
    0.82,
    0.828,
    0.831,
    0.834,
    0.838,
    0.834,
    0.835,
    0.837,
    0.846
)

# Format the first data cell explicitly (percentage with the same thin
# border used elsewhere in the sheet), then fan that formatting out to the
# rest of the column via copy/paste so every row ends up on one shared
# style record.
$ws.Cells.Item(2, 4).Value = $loadFactors[0]
$ws.Cells.Item(2, 4).NumberFormat = "0.00%"
$ws.Cells.Item(2, 4).Borders.Color = 0
$ws.Cells.Item(2, 4).Borders.LineStyle = 1

$ws.Range("D2").Copy()
for ($i = 1; $i -lt $loadFactors.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $loadFactors[$i]
    $ws.Cells.Item($row, 4).PasteSpecial(-4122)   # xlPasteFormats
}

# --- Restore the selection to what the authored file shows ---
$ws.Range("C3").Select()

# --- Printer settings added alongside the new column ---
$ws.PageSetup.Orientation = 1      # xlPortrait
$ws.PageSetup.PrintQuality = 300
